# petty-cashBook-2021.xlsx - "Update 29-Mei-2021, midday update."
# Adds daily petty-cash entries for 27/28/29-May-2021 (rows 22-42) and
# amends a few existing formulas in rows 19-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Amend existing rows 19-21
# ---------------------------------------------------------------------
$ws.Range("D19").Formula = "=60000+240000"
$ws.Range("C20").Formula = "=40965000+8892000+21622500+18858000"
$ws.Range("D21").Formula = "=40965000+21082500+9432500+640000+1100000+175000"

# ---------------------------------------------------------------------
# 2. New entries - 27 May 2021 (date serial 44343)
# ---------------------------------------------------------------------
# Row 22
$ws.Range("B22").Value = "SALES - cash/retail"
$ws.Range("C22").Formula = "=1212975+21250525-18858000"

# Row 23
$ws.Range("B23").Value = "SETOR KE BANK"
$ws.Range("D23").Formula = "=20000000"

# Row 24
$ws.Range("A24").Value = 44343
$ws.Range("B24").Value = "Wages Expense"
$ws.Range("D24").Formula = "=60000+280000"

# Row 25
$ws.Range("B25").Value = "A/R"
$ws.Range("C25").Formula = "=5000000+1750000+15900000+3000000+10932500"

# Row 26
$ws.Range("B26").Value = "TRANSFER BCA"
$ws.Range("D26").Formula = "=144000+1750000+4281500+401000+18900000"

# Row 27
$ws.Range("B27").Value = "PRIVE - andreas"
$ws.Range("D27").Value = 10000000

# Row 28
$ws.Range("B28").Value = "SALES - cash/retail"
$ws.Range("C28").Formula = "=12144525+8875975-10932500"

# Row 29
$ws.Range("B29").Value = "SELISIH - kurang"
$ws.Range("D29").Value = 290000

# Row 30
$ws.Range("B30").Value = "SETOR KE BANK"
$ws.Range("D30").Value = 11000000

# ---------------------------------------------------------------------
# 3. New entries - 28 May 2021 (date serial 44344)
# ---------------------------------------------------------------------
# Row 31
$ws.Range("A31").Value = 44344
$ws.Range("B31").Value = "Wages Expense"
$ws.Range("D31").Formula = "=60000+240000"

# Row 32
$ws.Range("B32").Value = "A/R"
$ws.Range("C32").Formula = "=10000000+10364500+143000+14124000"

# Row 33
$ws.Range("B33").Value = "TRANSFER BCA"
$ws.Range("D33").Formula = "=20507500+581000+135000"

# Row 34
$ws.Range("B34").Value = "FREIGHT OUT"
$ws.Range("D34").Formula = "=14400"

# Row 35
$ws.Range("B35").Value = "SALES - cash/retail"
$ws.Range("C35").Formula = "=23276125-500125-14124000"

# Row 36
$ws.Range("B36").Value = "SELISIH - lebih"
$ws.Range("C36").Value = 10000

# Row 37
$ws.Range("B37").Value = "SETOR KE BANK"
$ws.Range("D37").Formula = "=22000000"

# ---------------------------------------------------------------------
# 4. New entries - 29 May 2021 (date serial 44345)
# ---------------------------------------------------------------------
# Row 38
$ws.Range("A38").Value = 44345
$ws.Range("B38").Value = "Wages Expense"
$ws.Range("D38").Formula = "=60000"

# Row 39
$ws.Range("B39").Value = "SOLAR - KIJANG "
$ws.Range("D39").Formula = "=300000"

# Row 40
$ws.Range("B40").Value = "TRANSFER BCA"
$ws.Range("D40").Formula = "=50000000+8330000+561500+875000"

# Row 41
$ws.Range("B41").Value = "CHEQUE RECEIVED"
$ws.Range("D41").Formula = "=4043000"

# Row 42
$ws.Range("B42").Value = "A/R"
$ws.Range("C42").Formula = "=500000"

# ---------------------------------------------------------------------
# 5. Update the frozen-pane view / current selection
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("D58").Select()
